$wb = $excel.ActiveWorkbook

# "Test Steps" sheet: rows 10 and 14, column H (Results) go from PASS to FAIL
$wsSteps = $wb.Worksheets.Item("Test Steps")
$wsSteps.Range("H10").Value = "FAIL"
$wsSteps.Range("H14").Value = "FAIL"

# "Test Cases" sheet: rows 4 and 5, column D (Results) go from PASS to FAIL
$wsCases = $wb.Worksheets.Item("Test Cases")
$wsCases.Range("D4").Value = "FAIL"
$wsCases.Range("D5").Value = "FAIL"
